# edit.ps1 -- applies the "worked on stuff in class" commit to before.docx
#
# Summary of the change (see xml_diff):
#  1. Three existing bullets get their single run split into multiple runs
#     wrapped with <w:proofErr> spell-check markers (text is unchanged):
#       - "clean up year colums (2)"
#       - "create 2 dataframes"
#       - "check to make sure all events appear across the 15 olympics"
#  2. The "Individual sports vs teams column" bullet gets split into runs
#     with a <w:proofErr w:type="gramStart/gramEnd"/> wrap around "teams",
#     and the _GoBack bookmark that used to live in that paragraph moves
#     away from it.
#  3. The whole "Saturday" to-do sub-list is reordered/extended: the bullets
#     "Cleanup event column" ... "Merge" move up to sit right after
#     "Individual sports vs teams column", and two brand new bullets
#     ("Find and remove mixed events" / "Remove first two words of
#     subevents" -- the second one now carries the _GoBack bookmark) plus
#     a new empty bullet are appended after "Merge".
#
# Strategy: Word's Range.InsertXML, when called on a Range whose
# Start/End exactly match one or more whole paragraphs, replaces those
# paragraph(s) in place with the supplied OOXML (it does NOT just insert
# at the end of a paragraph the way a sub-range InsertXML does). That
# lets us rebuild each affected paragraph (or contiguous run of
# paragraphs) with exactly the run/proofErr/bookmark structure the diff
# calls for, while leaving every other paragraph byte-for-byte alone.

$d = $word.ActiveDocument

$pkgOpen = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage' xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document><w:body>"
$pkgClose = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

function Set-ParagraphXML($paragraph, [string]$innerXml) {
    $full = $paragraph.Range
    $full.InsertXML($pkgOpen + $innerXml + $pkgClose)
}

function Set-RangeXML($range, [string]$innerXml) {
    $range.InsertXML($pkgOpen + $innerXml + $pkgClose)
}

# ---------------------------------------------------------------------
# 1. "clean up year colums (2) - done!" -> split "colums" out with
#    spell-check proofErr wrapper; text unchanged.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "clean up year colums*") { $target = $p; break }
}
Set-ParagraphXML $target (
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:highlight w:val="green"/><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">clean up year </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>colums</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> (2)</w:t></w:r>' +
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> &#8211; done!</w:t></w:r>' +
    '</w:p>'
)

# ---------------------------------------------------------------------
# 2. "create 2 dataframes by sex" -> split "dataframes" out with
#    spell-check proofErr wrapper; text unchanged.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "create 2 dataframes*") { $target = $p; break }
}
Set-ParagraphXML $target (
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">create 2 </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>dataframes</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> by sex</w:t></w:r>' +
    '</w:p>'
)

# ---------------------------------------------------------------------
# 3. "check to make sure all events appear across the 15 olympics" ->
#    split "olympics" out with spell-check proofErr wrapper; text
#    unchanged.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "check to make sure all events*") { $target = $p; break }
}
Set-ParagraphXML $target (
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">check to make sure all events appear across the 15 </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>olympics</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
)

# ---------------------------------------------------------------------
# 4. The "Individual sports vs teams column" ... "Merge" block: split the
#    first bullet's run (grammar proofErr wrap around "teams", bookmark
#    removed from here), reorder the "Cleanup event column" .. "Merge"
#    bullets to sit right after it (unchanged), then append the two new
#    highlighted bullets (the bookmark now lives on the 2nd one) plus a
#    new empty bullet.
# ---------------------------------------------------------------------
$firstPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Individual sports vs teams column*") { $firstPara = $p; break }
}
$lastPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Merge`r") { $lastPara = $p; break }
}

$blockRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)

$newBlock = (
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Individual sports vs </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>teams</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> column</w:t></w:r>' +
    '</w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Cleanup event column</w:t></w:r>' +
    '</w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Import </w:t></w:r>' +
    '<w:r><w:t>world height data</w:t></w:r>' +
    '</w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>C</w:t></w:r>' +
    '<w:r><w:t>l</w:t></w:r>' +
    '<w:r><w:t>ean world data</w:t></w:r>' +
    '</w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Finish </w:t></w:r>' +
    '<w:r><w:t>Olympic cleaning</w:t></w:r>' +
    '</w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Merge</w:t></w:r>' +
    '</w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Find and remove mixed events</w:t></w:r>' +
    '</w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Remove first two words of subevents</w:t></w:r>' +
    '</w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr></w:p>'
)

Set-RangeXML $blockRange $newBlock

Write-Output "done"
